# Update AVTA course Excel file: add new course/promotion rows to the "courses" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# --- Row 2: Diploma of Civil Construction Design ---
$ws.Range("A2").Value = "RII50520"
$ws.Range("B2").Value = "111827M"
$ws.Range("C2").Value = "CIVIL CONSTRUCTION DESIGN"
$ws.Range("D2").Value = "DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("E2").Value = 52
$ws.Range("H2").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I2").Value = 10200
$ws.Range("J2").Value = "10,000 tuition fee + 200 handling fee"
$ws.Range("M2").Value = "TAS"
$ws.Range("R2").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "

# --- Row 3: Advanced Diploma of Civil Construction Design ---
$ws.Range("A3").Value = "RII60520"
$ws.Range("B3").Value = "111826A"
$ws.Range("C3").Value = "CIVIL CONSTRUCTION DESIGN"
$ws.Range("D3").Value = "ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("E3").Value = 104
$ws.Range("H3").Value = "88 wks tuition + 16 wks break"
$ws.Range("I3").Value = 20200
$ws.Range("J3").Value = "20,000 tuition fee + 200 handling fee"
$ws.Range("M3").Value = "TAS"
$ws.Range("R3").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "

# --- Row 4: Advanced Diploma of Information Technology (Telecommunications Network Engineering) ---
$ws.Range("A4").Value = "ICT60220"
$ws.Range("B4").Value = "111825B"
$ws.Range("C4").Value = "INFORMATION TECHNOLOGY"
$ws.Range("D4").Value = "ADVANCED DIPLOMA OF INFORMATION TECHNOLOGY `n(TELECOMMUNICATIONS NETWORK ENGINEERING) "
$ws.Range("E4").Value = 104
$ws.Range("H4").Value = "88 wks tuition + 16 wks break"
$ws.Range("I4").Value = 16200
$ws.Range("J4").Value = "16,000 tuition fee + 200 handling fee"
$ws.Range("M4").Value = "TAS"
$ws.Range("R4").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "

# --- Row 5: Package - Diploma + Advanced Diploma of Civil Construction Design ---
$ws.Range("A5").Value = "RII50520/RII60520"
$ws.Range("B5").Value = "111827M/111826A"
$ws.Range("C5").Value = "PACKAGES"
$ws.Range("D5").Value = "DIPLOMA OF CIVIL CONSTRUCTION DESIGN + ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("E5").Value = 104
$ws.Range("H5").Value = "88 wks tuition + 16 wks break"
$ws.Range("I5").Value = 20200
$ws.Range("J5").Value = "20,000 tuition fee + 200 handling fee"
$ws.Range("M5").Value = "TAS"
$ws.Range("R5").Value = "PROMOTIONS VALID UNTIL 28TH FEBRUARY 2023 "

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 45

# --- Number format for tuition column (I) = #,##0 (no wrap) ---
$ws.Range("I2:I5").NumberFormat = "#,##0"

# --- Number format + wrap text for tuitionDetail column (J) ---
$ws.Range("J2:J5").NumberFormat = "#,##0"
$ws.Range("J2:J5").WrapText = $true

# --- Wrap text (General number format) for other detail cells ---
$ws.Range("H2:H5").WrapText = $true
$ws.Range("R2:R5").WrapText = $true
$ws.Range("D4").WrapText = $true
$ws.Range("A5").WrapText = $true
$ws.Range("B5").WrapText = $true
$ws.Range("D5").WrapText = $true

# --- View settings: scroll so column F is the leftmost visible column, and select R14 ---
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("R14").Select()

$wb.Save()
